$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Update existing cells (row 91) ----
$ws.Cells.Item(91, 6).Value = 57144232.95   # F91
$ws.Cells.Item(91, 8).Value = 560           # H91

# ---- Update existing cells (row 92) ----
$ws.Cells.Item(92, 6).Value = 105460251.95  # F92
$ws.Cells.Item(92, 7).Value = 58            # G92
$ws.Cells.Item(92, 8).Value = 1259          # H92

# ---- Update existing cells (row 584) ----
$ws.Cells.Item(584, 2).Value = 27380          # B584
$ws.Cells.Item(584, 4).Value = 27450          # D584
$ws.Cells.Item(584, 6).Value = 3319949177.1   # F584
$ws.Cells.Item(584, 7).Value = 122410         # G584
$ws.Cells.Item(584, 8).Value = 3415           # H584

# ---- Helper to write a new data row, keeping column A as plain text ----
function Set-QuoteRow($Row, $Date, $Ultimo, $Apertura, $Maximo, $Minimo, $Monto, $Volumen, $Cantidad) {
    $cellA = $ws.Cells.Item($Row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $Date
    $cellA.ClearFormats()

    $ws.Cells.Item($Row, 2).Value = $Ultimo
    $ws.Cells.Item($Row, 3).Value = $Apertura
    $ws.Cells.Item($Row, 4).Value = $Maximo
    $ws.Cells.Item($Row, 5).Value = $Minimo
    $ws.Cells.Item($Row, 6).Value = $Monto
    $ws.Cells.Item($Row, 7).Value = $Volumen
    $ws.Cells.Item($Row, 8).Value = $Cantidad
}

# ---- New rows 586-590 ----
Set-QuoteRow 586 "2024-05-29" 27290       27601 27810       26810       6950046892       255707 4229
Set-QuoteRow 587 "2024-05-30" 29248       27619 29317       27619       7739473813.8     270208 6733
Set-QuoteRow 588 "2024-05-31" 29055       29250 29250       28000       10633833441      368237 6692
Set-QuoteRow 589 "2024-06-03" 28480       29000 29381.75    28300       3197527878.6     9      3401
Set-QuoteRow 590 "2024-06-04" 26739.25    27680 27850       26597.25    20875878603.75   772918 13330

$ws.Range("A1").Select() | Out-Null
